$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark numeric-looking price cells as Text so literal formatting (trailing zeros, etc.) is preserved
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.730.99'
$ws.Range("E2").Value = '  -4.19%  '
$ws.Range("D3").Value = '1.817.45'
$ws.Range("E3").Value = '  -3.05%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '277.56'
$ws.Range("E5").Value = '  -7.93%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '0.5091'
$ws.Range("E7").Value = '  -5.02%  '
$ws.Range("E8").Value = '  -6.07%  '
$ws.Range("E9").Value = '  -2.65%  '
$ws.Range("D10").Value = '0.06669'
$ws.Range("E10").Value = '  -7.21%  '
$ws.Range("D11").Value = '20.06'
$ws.Range("E11").Value = '  -7.01%  '
$ws.Range("D12").Value = '0.8282'
$ws.Range("E12").Value = '  -6.89%  '
$ws.Range("D13").Value = '0.07879'
$ws.Range("E13").Value = '  -3.27%  '
$ws.Range("D14").Value = '1.821.95'
$ws.Range("E14").Value = '  -3.48%  '
$ws.Range("D15").Value = '5.072'
$ws.Range("E15").Value = '  -4.66%  '
$ws.Range("D16").Value = '87.59'
$ws.Range("E16").Value = '  -6.32%  '
$ws.Range("D17").Value = '0.9998'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '14.11'
$ws.Range("E18").Value = '  -5.11%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000008030'
$ws.Range("E19").Value = '  -6.02%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = '25.770.24'
$ws.Range("E21").Value = '  -4.17%  '
$ws.Range("D22").Value = '4.741'
$ws.Range("E22").Value = '  -4.93%  '
$ws.Range("D23").Value = '10.01'
$ws.Range("E23").Value = '  -5.93%  '
$ws.Range("E24").Value = '  -4.86%  '
$ws.Range("E25").Value = '  -2.66%  '
$ws.Range("D26").Value = '2.201'
$ws.Range("E26").Value = '  -4.33%  '
$ws.Range("D27").Value = '1.674'
$ws.Range("E27").Value = '  -3.32%  '
$ws.Range("E28").Value = '  -5.37%  '
$ws.Range("D29").Value = '109.53'
$ws.Range("E29").Value = '  -3.96%  '
$ws.Range("D30").Value = '4.333'
$ws.Range("E30").Value = '  -8.26%  '
$ws.Range("D31").Value = '4.234'
$ws.Range("E31").Value = '  -8.29%  '
$ws.Range("D32").Value = '0.08779'
$ws.Range("E32").Value = '  -4.02%  '
$ws.Range("D33").Value = '0.04889'
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("D34").Value = '0.7277'
$ws.Range("E34").Value = '  -10.60%  '
$ws.Range("D35").Value = '1.138'
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("D36").Value = '2.868'
$ws.Range("E36").Value = '  -2.68%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '3.124'
$ws.Range("E37").Value = '  -2.83%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.373'
$ws.Range("E38").Value = '  -9.15%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01853'
$ws.Range("E39").Value = '  -5.19%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5177'
$ws.Range("E40").Value = '  -14.29%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.9642'
$ws.Range("E41").Value = '  -9.94%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.218'
$ws.Range("E42").Value = '  -6.28%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '110.33'
$ws.Range("E43").Value = '  -4.20%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '8.030'
$ws.Range("E44").Value = '  -10.09%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.4553'
$ws.Range("E46").Value = '  -10.51%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1365'
$ws.Range("E47").Value = '  -8.63%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '36.50'
$ws.Range("E48").Value = '  -3.22%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.259'
$ws.Range("E49").Value = '  -6.88%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.503'
$ws.Range("E50").Value = '  -8.07%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05840'
$ws.Range("E51").Value = '  -3.49%  '

# Restore default (Normal) style on the cells we temporarily marked as Text,
# so only the cell VALUE changes relative to the original (style index stays 0).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
